# Generate Report for Handback
# Adds a new handed-back file (614fb894-1524-4b63-a580-394b1a45dc50.md) as
# row 4 to the "Overview", "zh-cn" and "de-de" sheets, mirroring the
# existing rows for the already-in-sync file (1ce831e9-...).

$wb = $excel.ActiveWorkbook

$fileName   = "614fb894-1524-4b63-a580-394b1a45dc50.md"
$pathName   = "e2e\614fb894-1524-4b63-a580-394b1a45dc50.md"
$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G4").Value = "2017-02-17 09:36:07"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/614fb894-1524-4b63-a580-394b1a45dc50.md",
    "",
    "",
    $pathName
) | Out-Null
$wsOverview.Range("B4").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = $statusInSync
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "614fb894-1524-4b63-a580-394b1a45dc50.4fe3a609b5fe7955a0db8d224017350ac6d1704d.zh-cn.xlf"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H4").Value = "2017-02-17 09:35:50"
$wsZhCn.Range("I4").Value = "'"
$wsZhCn.Range("K4").Value = "614fb894-1524-4b63-a580-394b1a45dc50.4fe3a609b5fe7955a0db8d224017350ac6d1704d.zh-cn.xlf"
$wsZhCn.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value = "2017-02-17 09:36:44"
$wsZhCn.Range("M4").Value = "'"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'True"
$wsZhCn.Range("P4").Value = "'"
$wsZhCn.Range("Q4").Value = "'False"
$wsZhCn.Range("R4").Value = "'"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/614fb894-1524-4b63-a580-394b1a45dc50.md",
    "",
    "",
    $fileName
) | Out-Null
$wsZhCn.Range("A4").Style = "Hyperlink"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("J4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/0000000000000000000000000000000000000000/e2e/614fb894-1524-4b63-a580-394b1a45dc50.md",
    "",
    "",
    $fileName
) | Out-Null
$wsZhCn.Range("J4").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = $statusInSync
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "614fb894-1524-4b63-a580-394b1a45dc50.4fe3a609b5fe7955a0db8d224017350ac6d1704d.de-de.xlf"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H4").Value = "2017-02-17 09:36:07"
$wsDeDe.Range("I4").Value = "'"
$wsDeDe.Range("K4").Value = "614fb894-1524-4b63-a580-394b1a45dc50.4fe3a609b5fe7955a0db8d224017350ac6d1704d.de-de.xlf"
$wsDeDe.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value = "2017-02-17 09:37:08"
$wsDeDe.Range("M4").Value = "'"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'True"
$wsDeDe.Range("P4").Value = "'"
$wsDeDe.Range("Q4").Value = "'False"
$wsDeDe.Range("R4").Value = "'"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/614fb894-1524-4b63-a580-394b1a45dc50.md",
    "",
    "",
    $fileName
) | Out-Null
$wsDeDe.Range("A4").Style = "Hyperlink"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("J4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/0000000000000000000000000000000000000000/e2e/614fb894-1524-4b63-a580-394b1a45dc50.md",
    "",
    "",
    $fileName
) | Out-Null
$wsDeDe.Range("J4").Style = "Hyperlink"
